$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 349, shifting existing rows 349:435 down to 350:436
$ws.Rows("349:349").Insert()

# Populate the new row 349 with the new data record
$ws.Range("A349").Value = 3
$ws.Range("B349").Value = "Femacal de La Calera"
$ws.Range("C349").Value = "Coquimbo"
$ws.Range("D349").Value = 44855
$ws.Range("E349").Value = 5
$ws.Range("F349").Value = 100112043
$ws.Range("G349").Value = "Pepino ensalada"
$ws.Range("H349").Value = "Sin especificar"
$ws.Range("I349").Value = "Primera"
$ws.Range("J349").Value = 125
$ws.Range("K349").Value = 17000
$ws.Range("L349").Value = 17500
$ws.Range("M349").Value = 17260
$ws.Range("N349").Value = '$/caja 60 unidades'
$ws.Range("O349").Value = "Región de Arica y Parinacota"
$ws.Range("P349").Value = 288
$ws.Range("Q349").Value = 60
$ws.Range("R349").Value = "Hortaliza"
